$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.818.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.634.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0779'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("E12").Value = '  -0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.639.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.860.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.826.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("E22").Value = '  +1.12%  '

$ws.Range("E23").Value = '  +2.86%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("E25").Value = '  -2.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.05%  '

$ws.Range("E27").Value = '  -3.70%  '

$ws.Range("E28").Value = '  +1.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("E30").Value = '  +0.31%  '

$ws.Range("E31").Value = '  +1.54%  '

$ws.Range("E32").Value = '  +1.36%  '

$ws.Range("E33").Value = '  +1.80%  '

$ws.Range("E34").Value = '  +1.50%  '

$ws.Range("E35").Value = '  +0.31%  '

$ws.Range("E36").Value = '  +0.84%  '

$ws.Range("E37").Value = '  +0.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.550'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.115.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("E40").Value = '  +0.60%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.799'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.91%  '

$ws.Range("E48").Value = '  -5.30%  '

$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("E51").Value = '  +0.14%  '
